$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 8969  # F2: 8960 -> 8969
$ws.Cells.Item(4, 6).Value = 6601  # F4: 6597 -> 6601
$ws.Cells.Item(5, 6).Value = 174  # F5: 173 -> 174
$ws.Cells.Item(6, 6).Value = 2129  # F6: 2127 -> 2129
$ws.Cells.Item(7, 6).Value = 596  # F7: 595 -> 596
$ws.Cells.Item(13, 6).Value = 1  # F13: 0 -> 1
$ws.Cells.Item(16, 6).Value = 8862  # F16: 8845 -> 8862
$ws.Cells.Item(20, 6).Value = 118  # F20: 117 -> 118
$ws.Cells.Item(25, 6).Value = 81  # F25: 80 -> 81
$ws.Cells.Item(28, 6).Value = 1033  # F28: 1032 -> 1033
$ws.Cells.Item(29, 6).Value = 9  # F29: 8 -> 9
$ws.Cells.Item(30, 6).Value = 67  # F30: 65 -> 67
$ws.Cells.Item(31, 6).Value = 549  # F31: 548 -> 549
$ws.Cells.Item(33, 6).Value = 18  # F33: 17 -> 18
$ws.Cells.Item(34, 6).Value = 538  # F34: 537 -> 538
$ws.Cells.Item(35, 6).Value = 2299  # F35: 2290 -> 2299
$ws.Cells.Item(36, 6).Value = 873  # F36: 872 -> 873
$ws.Cells.Item(37, 6).Value = 536  # F37: 531 -> 536
$ws.Cells.Item(41, 6).Value = 285  # F41: 281 -> 285
$ws.Cells.Item(44, 6).Value = 1047  # F44: 1044 -> 1047
$ws.Cells.Item(45, 6).Value = 97  # F45: 95 -> 97
$ws.Cells.Item(46, 6).Value = 99  # F46: 98 -> 99
$ws.Cells.Item(47, 6).Value = 15  # F47: 11 -> 15
$ws.Cells.Item(48, 6).Value = 77  # F48: 76 -> 77

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(14, 6).Value = 10  # F14: 9 -> 10

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 329  # F4: 328 -> 329

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 8969  # F3: 8960 -> 8969
$ws.Cells.Item(4, 6).Value = 329  # F4: 328 -> 329
$ws.Cells.Item(6, 6).Value = 6601  # F6: 6597 -> 6601
$ws.Cells.Item(7, 6).Value = 174  # F7: 173 -> 174
$ws.Cells.Item(8, 6).Value = 2129  # F8: 2127 -> 2129
$ws.Cells.Item(11, 6).Value = 596  # F11: 595 -> 596
$ws.Cells.Item(19, 6).Value = 8862  # F19: 8845 -> 8862
$ws.Cells.Item(26, 6).Value = 81  # F26: 80 -> 81
$ws.Cells.Item(28, 6).Value = 1033  # F28: 1032 -> 1033
$ws.Cells.Item(29, 6).Value = 9  # F29: 8 -> 9
$ws.Cells.Item(30, 6).Value = 67  # F30: 65 -> 67
$ws.Cells.Item(32, 6).Value = 549  # F32: 548 -> 549
$ws.Cells.Item(34, 6).Value = 18  # F34: 17 -> 18
$ws.Cells.Item(35, 6).Value = 538  # F35: 537 -> 538
$ws.Cells.Item(36, 6).Value = 2299  # F36: 2290 -> 2299
$ws.Cells.Item(37, 6).Value = 873  # F37: 872 -> 873
$ws.Cells.Item(38, 6).Value = 10  # F38: 9 -> 10
$ws.Cells.Item(40, 6).Value = 536  # F40: 531 -> 536
$ws.Cells.Item(41, 6).Value = 285  # F41: 281 -> 285
$ws.Cells.Item(43, 6).Value = 77  # F43: 76 -> 77
